$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells keep their text formatting so values
# like "595.07" or "0.0000234" are not auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.067.11"
$ws.Range("E2").Value = "  +3.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.031.58"
$ws.Range("E3").Value = "  +1.71%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.07"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.65"
$ws.Range("E6").Value = "  +7.45%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.030.25"
$ws.Range("E8").Value = "  +1.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.97"
$ws.Range("E10").Value = "  +16.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  +2.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  +2.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").Value = "  +3.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.59"
$ws.Range("E14").Value = "  +4.33%  "

$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.539.94"
$ws.Range("E16").Value = "  +2.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("E17").Value = "  +3.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.959.42"
$ws.Range("E18").Value = "  +2.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.032.90"
$ws.Range("E19").Value = "  +1.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.75"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("E21").Value = "  +2.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("E22").Value = "  +2.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  +3.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  +7.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.40"
$ws.Range("E27").Value = "  +4.02%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("E29").Value = "  +6.23%  "

$ws.Range("E30").Value = "  +11.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.71"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.65"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("E34").Value = "  +3.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0863"
$ws.Range("E35").Value = "  +6.77%  "

$ws.Range("E36").Value = "  +3.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.91"
$ws.Range("E37").Value = "  +2.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.15"
$ws.Range("E38").Value = "  +11.24%  "

$ws.Range("E39").Value = "  +9.91%  "

$ws.Range("E40").Value = "  +3.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.50"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.11"
$ws.Range("E42").Value = "  +1.71%  "

$ws.Range("E43").Value = "  +15.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.50"
$ws.Range("E44").Value = "  +15.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.29"
$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0361"
$ws.Range("E46").Value = "  +3.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.722.65"
$ws.Range("E47").Value = "  +0.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.77"
$ws.Range("E48").Value = "  +2.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.83"
$ws.Range("E49").Value = "  +11.36%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.29"
$ws.Range("E51").Value = "  +7.56%  "

# Rows 24 and 25 swap coins (RenderToken <-> Litecoin) with updated data
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.32"
$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.47"
$ws.Range("E25").Value = "  +9.53%  "
